$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "ID" column becomes "NO" ---
$ws.Range("A1").Value = "NO"

# --- Data fix: Eselon for row 3 (Bebek) becomes a literal number instead of text "II" ---
$ws.Range("H3").Value = 2

# --- Insert two new columns (Berangkat / Kembali) before the old "Nama Rekening" column ---
$ws.Columns("J:K").Insert()

# Headers for the new columns (inherit the bold header style already propagated by Insert)
$ws.Range("J1").Value = "Berangkat"
$ws.Range("K1").Value = "Kembali"

# Data rows for the new columns - departure / return dates
$ws.Range("J2").Value = 45586
$ws.Range("K2").Value = 45649
$ws.Range("J3").Value = 45587
$ws.Range("K3").Value = 45650
$ws.Range("J4").Value = 45588
$ws.Range("K4").Value = 45651
$ws.Range("J5").Value = 45589
$ws.Range("K5").Value = 45652
$ws.Range("J6").Value = 45590
$ws.Range("K6").Value = 45653

# Apply date formatting to the new header + data cells
$ws.Range("J1:K1").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("J2:K6").NumberFormat = "yyyy\-mm\-dd;@"

# Column widths for the new columns (~14.33 chars, matching the sibling date columns)
$ws.Columns("J:K").ColumnWidth = 13.498697916666666

# --- Fix up the Golongan/Ruang reference list validation formula, which referenced
#     column P before the column insert shifted that data to column R ---
$ws.Range("F1:F1048576").Validation.Formula1 = "=`$R`$2:`$R`$18"

Write-Output "done"
